# This workbook tracks daily "Zapallo italiano" (zucchini) price records for
# "Macroferia Regional de Talca". The commit adds two new daily records into
# the middle of the existing data block:
#   - a new record (dated 2022-01-08, serial 44567) inserted as row 165
#   - a new record (dated 2022-01-09, serial 44568) inserted as row 256
# Every existing row at/after each insertion point shifts down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert the first new row at position 165 -----------------------------
$ws.Rows.Item(165).Insert()

$ws.Range("A165").Value = 5
$ws.Range("B165").Value = "Macroferia Regional de Talca"
$ws.Range("C165").Value = "Maule"
$ws.Range("D165").Value = 44567
$ws.Range("E165").Value = 7
$ws.Range("F165").Value = 100112032
$ws.Range("G165").Value = "Zapallo italiano"
$ws.Range("H165").Value = "Sin especificar"
$ws.Range("I165").Value = "Primera"
$ws.Range("J165").Value = 300
$ws.Range("K165").Value = 5000
$ws.Range("L165").Value = 5000
$ws.Range("M165").Value = 5000
$ws.Range("N165").Value = "`$/caja 60 unidades"
$ws.Range("O165").Value = "Región del Maule"
$ws.Range("P165").Value = 83
$ws.Range("Q165").Value = 60
$ws.Range("R165").Value = "Hortaliza"

# --- Insert the second new row at position 256 (post first-insert index) --
$ws.Rows.Item(256).Insert()

$ws.Range("A256").Value = 5
$ws.Range("B256").Value = "Macroferia Regional de Talca"
$ws.Range("C256").Value = "Maule"
$ws.Range("D256").Value = 44568
$ws.Range("E256").Value = 7
$ws.Range("F256").Value = 100112032
$ws.Range("G256").Value = "Zapallo italiano"
$ws.Range("H256").Value = "Sin especificar"
$ws.Range("I256").Value = "Primera"
$ws.Range("J256").Value = 300
$ws.Range("K256").Value = 6000
$ws.Range("L256").Value = 6000
$ws.Range("M256").Value = 6000
$ws.Range("N256").Value = "`$/caja 60 unidades"
$ws.Range("O256").Value = "Región del Maule"
$ws.Range("P256").Value = 100
$ws.Range("Q256").Value = 60
$ws.Range("R256").Value = "Hortaliza"
